# Applies a cyclic re-shuffle of the per-row data (columns D, L, M, N, O, P, Q, R, S, T)
# across rows 2, 3, 5, 6, 7, 8, 9, 10, 11, 12 of Sheet1, matching the target diff.
# Rows 4 and 13 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> index map for the columns that actually change.
# D=4, L=12, M=13, N=14, O=15, P=16, Q=17, R=18, S=19, T=20

function Get-RowData($row) {
    return @{
        D = $ws.Cells.Item($row, 4).Value2
        L = $ws.Cells.Item($row, 12).Value2
        M = $ws.Cells.Item($row, 13).Value2
        N = $ws.Cells.Item($row, 14).Value2
        O = $ws.Cells.Item($row, 15).Value2
        P = $ws.Cells.Item($row, 16).Value2
        Q = $ws.Cells.Item($row, 17).Value2
        R = $ws.Cells.Item($row, 18).Value2
        S = $ws.Cells.Item($row, 19).Value2
        T = $ws.Cells.Item($row, 20).Value2
    }
}

function Set-RowData($row, $data) {
    $ws.Cells.Item($row, 4).Value2 = $data.D
    $ws.Cells.Item($row, 12).Value2 = $data.L
    $ws.Cells.Item($row, 13).Value2 = $data.M
    $ws.Cells.Item($row, 14).Value2 = $data.N
    $ws.Cells.Item($row, 15).Value2 = $data.O
    $ws.Cells.Item($row, 16).Value2 = $data.P
    $ws.Cells.Item($row, 17).Value2 = $data.Q
    $ws.Cells.Item($row, 18).Value2 = $data.R
    $ws.Cells.Item($row, 19).Value2 = $data.S
    $ws.Cells.Item($row, 20).Value2 = $data.T
}

# Snapshot all affected rows before mutating any of them.
$row2 = Get-RowData 2
$row3 = Get-RowData 3
$row5 = Get-RowData 5
$row6 = Get-RowData 6
$row7 = Get-RowData 7
$row8 = Get-RowData 8
$row9 = Get-RowData 9
$row10 = Get-RowData 10
$row11 = Get-RowData 11
$row12 = Get-RowData 12

# Cyclic mapping derived from the diff (old row -> new row):
#   2  -> 11, 11 -> 2        (swap)
#   3  -> 12, 12 -> 3        (swap)
#   5  -> 7,  7  -> 9, 9  -> 5   (3-cycle)
#   6  -> 8,  8  -> 10, 10 -> 6  (3-cycle)
Set-RowData 11 $row2
Set-RowData 2  $row11

Set-RowData 12 $row3
Set-RowData 3  $row12

Set-RowData 7 $row5
Set-RowData 9 $row7
Set-RowData 5 $row9

Set-RowData 8  $row6
Set-RowData 10 $row8
Set-RowData 6  $row10
